$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 / Row 21 swap: Uniswap <-> ShibaInu ---
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.60"
$ws.Range("E21").Value = "  -0.72%  "

# --- Row 42 / Row 43 swap: RenderToken <-> NEARProtocol ---
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("E42").Value = "  -2.09%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.87"
$ws.Range("E43").Value = "  -0.29%  "

# --- Remaining rows: Price (D) and Volume(1h) (E) updates ---

# Row 2
$ws.Range("D2").Value = "42.698.29"
$ws.Range("E2").Value = "  -0.43%  "

# Row 3
$ws.Range("D3").Value = "2.521.75"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.60"
$ws.Range("E5").Value = "  +4.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.75"
$ws.Range("E6").Value = "  -3.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.582"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("E9").Value = "  -1.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.27"
$ws.Range("E10").Value = "  -2.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12
$ws.Range("E12").Value = "  +1.46%  "

# Row 13
$ws.Range("E13").Value = "  -3.20%  "

# Row 14
$ws.Range("D14").Value = "2.908.40"
$ws.Range("E14").Value = "  -1.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.55"
$ws.Range("E15").Value = "  +3.25%  "

# Row 16
$ws.Range("D16").Value = "2.515.56"
$ws.Range("E16").Value = "  -1.70%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.861"
$ws.Range("E17").Value = "  -2.43%  "

# Row 18
$ws.Range("D18").Value = "42.686.39"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("E19").Value = "  -6.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.42"
$ws.Range("E22").Value = "  -0.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.86"
$ws.Range("E23").Value = "  -1.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.00"
$ws.Range("E24").Value = "  +0.42%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("E25").Value = "  -2.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.06"
$ws.Range("E26").Value = "  -2.82%  "

# Row 27
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  +13.17%  "

# Row 29
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.83"
$ws.Range("E30").Value = "  +1.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.92"
$ws.Range("E31").Value = "  -1.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.50"
$ws.Range("E32").Value = "  -1.86%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.62"
$ws.Range("E33").Value = "  +3.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.35"
$ws.Range("E34").Value = "  +0.64%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.10"
$ws.Range("E35").Value = "  -3.68%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0790"
$ws.Range("E36").Value = "  -2.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.61"
$ws.Range("E37").Value = "  -4.99%  "

# Row 38
$ws.Range("E38").Value = "  -1.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.30"
$ws.Range("E39").Value = "  -7.50%  "

# Row 40
$ws.Range("E40").Value = "  +0.53%  "

# Row 41
$ws.Range("E41").Value = "  -1.84%  "

# Row 44
$ws.Range("E44").Value = "  -1.01%  "

# Row 45
$ws.Range("E45").Value = "  +0.08%  "

# Row 46
$ws.Range("D46").Value = "2.041.91"
$ws.Range("E46").Value = "  -2.20%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.56"
$ws.Range("E47").Value = "  -2.36%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.81"
$ws.Range("E48").Value = "  -2.76%  "

# Row 49
$ws.Range("D49").Value = "2.764.59"
$ws.Range("E49").Value = "  -1.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.26"
$ws.Range("E50").Value = "  -1.83%  "

# Row 51
$ws.Range("E51").Value = "  +0.00%  "
